# Release Form - F5.xlsx
# Re-home the form from the "Change Request" (CR) series to the
# "Software Development Lifecycle" (SD) series:
#   - rename the worksheet tab F-SW-CR-05 -> F-SW-SD-05
#   - carry the Print_Area defined name over to the renamed sheet
#   - move the live selection to G49
#   - bump the footer's revision date stamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab name: F-SW-CR-05 -> F-SW-SD-05
$ws.Name = "F-SW-SD-05"

# _xlnm.Print_Area embeds the sheet name in its formula ('F-SW-CR-05'!$A$1:$E$13);
# re-asserting PrintArea rewrites that defined name against the new sheet name.
$ws.PageSetup.PrintArea = '$A$1:$E$13'

# Selection moves from H5 to G49.
$ws.Range("G49").Select() | Out-Null

# Footer revision date: (0/0/2025) -> (01/10/2025)
$ws.PageSetup.RightFooter = "&14Rev:0(01/10/2025)"
